# Update bond trade record:
# - Row 23 now represents a DATE(2014,8,1) cash-flow of -49397.11 (was DATE(2014,7,31) / 19550.24)
# - Row 24 becomes a new trade record: DATE(2014,8,11) cash-flow of 69227.75
# - Summary formulas in row 2/row 3 extend their ranges to include row 24 (and row 25 for XIRR on column B)
# - A new blank row 44 is appended below the existing data so the sheet grows by one row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23: change date and amount ---
$ws.Range("A23").Formula = "=DATE(2014,8,1)"
$ws.Range("F23").Value = -49397.11
# G23 keeps its existing formula =SUM(F23); value recalculates automatically

# --- Row 24: new trade record ---
$ws.Range("A24").Formula = "=DATE(2014,8,11)"
$ws.Range("A24").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("F24").Value = 69227.75
$ws.Range("G24").Formula = "=SUM(F24)"

# --- Row 2 summary formulas: extend ranges to row 24 ---
$ws.Range("F2").Formula = "=SUM(F22:F24)"
$ws.Range("G2").Formula = "=SUM(G4:G24)"

# --- Row 3 summary formulas: extend ranges ---
$ws.Range("B3").Formula = "=XIRR(B4:B25,A4:A25)"
$ws.Range("F3").Formula = "=XIRR(F22:F24,A22:A24)"
$ws.Range("G3").Formula = "=XIRR(G4:G24,A4:A24)"

# --- New blank row 44, matching the style/number format of the rows above it ---
$ws.Range("B44:G44").NumberFormat = "0.00"

# --- Selection moves to G2 (last-edited cell context) ---
$ws.Range("G2").Select()
